# Add handling errors on input:
# Insert a new "s1200" worksheet between "ACID VALUE" and "HF SERIES" that
# carries the header row for the new input-validation workflow (LOT, Step,
# Suhu, Berat Sample, Jumlah Titran, Faktor Buret, Faktor NaOH, AV,
# Instruksi, Operator QC), and make it the active sheet/tab.

$wb = $excel.ActiveWorkbook

$acidValue = $wb.Worksheets.Item("ACID VALUE")

# New sheet goes right after "ACID VALUE" (i.e. before "HF SERIES").
$newSheet = $wb.Worksheets.Add($null, $acidValue)
$newSheet.Name = "s1200"

# Fill in the header row. The fill order matters for shared-string layout
# (matches how the strings ended up appended to sharedStrings.xml): the
# "Step" (B1) and "Operator QC" (J1) headers were the two typed last, in
# that order.
$newSheet.Range("A1").Value = "LOT"
$newSheet.Range("C1").Value = "Suhu"
$newSheet.Range("D1").Value = "Berat Sample (gr)"
$newSheet.Range("E1").Value = "Jumlah Titran (mL)"
$newSheet.Range("F1").Value = "Faktor Buret"
$newSheet.Range("G1").Value = "Faktor NaOH"
$newSheet.Range("H1").Value = "AV"
$newSheet.Range("I1").Value = "Instruksi"
$newSheet.Range("B1").Value = "Step"
$newSheet.Range("J1").Value = "Operator QC"

# Auto-fit the columns that hold the longer headers, like Excel does after
# typing values into an empty column.
$newSheet.Range("D:D").AutoFit() | Out-Null
$newSheet.Range("E:E").AutoFit() | Out-Null
$newSheet.Range("F:F").AutoFit() | Out-Null
$newSheet.Range("G:G").AutoFit() | Out-Null
$newSheet.Range("H:H").AutoFit() | Out-Null
$newSheet.Range("J:J").AutoFit() | Out-Null

# Leave the selection parked where the author left it and make this new
# sheet the active tab.
$newSheet.Range("N11").Select() | Out-Null
$newSheet.Activate()
